$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.405.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.50%  "

# Row 3
$ws.Range("D3").Value = "'3.170.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "

# Row 5
$ws.Range("D5").Value = "'596.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.44%  "

# Row 6
$ws.Range("D6").Value = "'135.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "

# Row 7
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").Value = "'3.172.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.60%  "

# Row 9
$ws.Range("E9").Value = "  +1.67%  "

# Row 10
$ws.Range("E10").Value = "  -1.19%  "

# Row 11
$ws.Range("D11").Value = "'5.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "

# Row 12
$ws.Range("D12").Value = "'0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "

# Row 13
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.34%  "

# Row 14
$ws.Range("D14").Value = "'34.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.93%  "

# Row 15
$ws.Range("D15").Value = "'3.692.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "

# Row 16
$ws.Range("E16").Value = "  -0.10%  "

# Row 17
$ws.Range("D17").Value = "'3.172.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.83%  "

# Row 18
$ws.Range("D18").Value = "'63.396.81"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'6.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.00%  "

# Row 20
$ws.Range("D20").Value = "'461.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.46%  "

# Row 21
$ws.Range("D21").Value = "'13.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.81%  "

# Row 22
$ws.Range("D22").Value = "'0.696"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "

# Row 23
$ws.Range("D23").Value = "'7.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'83.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.01%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'13.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.35%  "

# Row 26
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

# Row 27
$ws.Range("D27").Value = "'2.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.07%  "

# Row 28
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "

# Row 29
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "'6.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.87%  "

# Row 30
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'2.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.41%  "

# Row 31
$ws.Range("D31").Value = "'7.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.45%  "

# Row 32
$ws.Range("D32").Value = "'27.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "

# Row 33
$ws.Range("E33").Value = "  -1.14%  "

# Row 34
$ws.Range("E34").Value = "  -1.08%  "

# Row 35
$ws.Range("E35").Value = "  -2.47%  "

# Row 36
$ws.Range("D36").Value = "'5.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.36%  "

# Row 37
$ws.Range("D37").Value = "'0.0₃0740"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.67%  "

# Row 38
$ws.Range("D38").Value = "'51.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.19%  "

# Row 39
$ws.Range("D39").Value = "'0.0390"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.88%  "

# Row 40
$ws.Range("D40").Value = "'8.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.17%  "

# Row 41
$ws.Range("E41").Value = "  -2.03%  "

# Row 42
$ws.Range("D42").Value = "'2.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.32%  "

# Row 43
$ws.Range("D43").Value = "'391.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.83%  "

# Row 44
$ws.Range("D44").Value = "'2.788.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.37%  "

# Row 45
$ws.Range("D45").Value = "'0.251"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.14%  "

# Row 46
$ws.Range("D46").Value = "'127.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.91%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "

# Row 48
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'35.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.83%  "

# Row 49
$ws.Range("D49").Value = "'2.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.99%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.112"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'25.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.58%  "
